$wb = $excel.ActiveWorkbook

# New handoff/handback timestamps generated for the "3026501a..." file
# during this handback report run. The other file (f1ebf36c...) was not
# touched in this run, so its timestamps stay as-is.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-18 04:49:12"
$wsZhCn.Range("K2").Value = "2016-08-18 04:49:32"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-18 04:49:17"
$wsDeDe.Range("K2").Value = "2016-08-18 04:49:39"

# Overview sheet mirrors the de-de "Correspond Handoff Datetime" as the
# "Latest HO Xliff Generate Date" for each source file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-18 04:49:17"
